$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet 2: Preallokering ---
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Preallokering"

$ws2.Range("A1").Value = "Medarbejder"
$ws2.Range("B1").Value = "The Change"
$ws2.Range("C1").Value = "LEG-DHC"
$ws2.Range("D1").Value = "LTDE-repBC"
$ws2.Range("E1").Value = "COOLGEOHEAT II"
$ws2.Range("F1").Value = "HEATCODE"

$ws2names = @(
  "Søren Erbs Poulsen (SOEB)",
  "Mathias Larsen (MATL)",
  "Rune Kier Nielsen (RUNI)",
  "Kristoffer Bested Nielsen (KRI)",
  "Søren Andersen (SSSA)",
  "Karl Woldum Tordrup (KART)",
  "Marton Major (MMAJ)",
  "NN (ufordelt)"
)
for ($i = 0; $i -lt $ws2names.Count; $i++) {
  $row = $i + 2
  $ws2.Range("A$row").Value = $ws2names[$i]
}

# data values (B:F), row2 is 50/50/0/0/0, others all 0
$ws2.Range("B2:F9").Value = 0
$ws2.Range("B2").Value = 50
$ws2.Range("C2").Value = 50

$ws2.Range("B2:F9").NumberFormat = "0"

# --- Sheet 3: Eksterne timer og porteføljer ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Eksterne timer og porteføljer"

$ws3.Range("A1").Value = "Medarbejder"
$ws3.Range("C1").Value = "Portefølje [t]"
$ws3.Range("B1").Value = "Ekstern [t]"

$ws3names = @(
  "Søren Erbs Poulsen (SOEB)",
  "Mathias Larsen (MATL)",
  "Rune Kier Nielsen (RUNI)",
  "Kristoffer Bested Nielsen (KRI)",
  "Søren Andersen (SSSA)",
  "Karl Woldum Tordrup (KART)",
  "Marton Major (MMAJ)",
  "NN (ufordelt)"
)
for ($i = 0; $i -lt $ws3names.Count; $i++) {
  $row = $i + 2
  $ws3.Range("A$row").Value = $ws3names[$i]
}

# B column: external hours
$ws3.Range("B2").Value = 50
$ws3.Range("B3").Value = 0
$ws3.Range("B4").Value = 480
$ws3.Range("B5").Value = 547
$ws3.Range("B6").Value = 500
$ws3.Range("B7").Value = 0
$ws3.Range("B8").Value = 0
$ws3.Range("B9").Value = 0

# C column: portfolio hours
$ws3.Range("C2").Formula = "=1315/2"
$ws3.Range("C3:C6").Formula = "=1315/2"
$ws3.Range("C7").Value = 100
$ws3.Range("C8").Formula = "=1315/2"
$ws3.Range("C9").Value = 10000

$ws3.Range("B2:C9").NumberFormat = "0"

# --- Column widths (approximate autofit) ---
$ws2.Columns.Item(1).ColumnWidth = 27.7109375
$ws2.Columns.Item(2).ColumnWidth = 11.140625
$ws2.Columns.Item(3).ColumnWidth = 8.5703125
$ws2.Columns.Item(4).ColumnWidth = 11
$ws2.Columns.Item(5).ColumnWidth = 15.7109375
$ws2.Columns.Item(6).ColumnWidth = 10.28515625

$ws3.Columns.Item(1).ColumnWidth = 27.7109375
$ws3.Columns.Item(2).ColumnWidth = 13.42578125
$ws3.Columns.Item(3).ColumnWidth = 12

# --- Selections to match final saved view state ---
[void]$ws1.Range("A1:F9").Select()
[void]$ws2.Range("A1:C9").Select()
[void]$ws3.Range("B1").Select()
[void]$ws3.Activate()


